# 8.10.1.xlsx "Add files via upload" re-edit:
#   - append a new "2023" column (T) mirroring the existing S (2022) column
#   - convert the now-historical shared-formula values in rows 4 and 5
#     (indicator calc = count / adult population * 100000) to plain
#     cached numbers, matching updated source figures
#   - refresh the adult-population row (row 8) figures for years 2010-2020
#   - a couple of row-height tweaks that came along with the re-save
#   - reset the sheet's saved cursor position to A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Add column T (2023) by cloning column S's formatting, then
#    overwrite with the real 2023 figures for every populated row.
# ---------------------------------------------------------------------
$ws.Range("S3:S8").Copy($ws.Range("T3:T8"))

$ws.Range("T3").Value = 2023
$ws.Range("T4").Value = 6.7904451646088795
$ws.Range("T5").Value = 47.957518975050206
$ws.Range("T6").Value = 320
$ws.Range("T7").Value = 2260
$ws.Range("T8").Value = 4712504

# ---------------------------------------------------------------------
# 2. Row 4 ("a) branches per 100,000 adults"): the cells D4:Q4 used to
#    hold a shared formula (=col6/col8*100000). They are now stored as
#    plain numbers with refreshed figures, so drop the formulas and
#    write the updated literal results.
# ---------------------------------------------------------------------
$ws.Range("D4").Value = 6.1074687240787666
$ws.Range("E4").Value = 6.5454292116044552
$ws.Range("F4").Value = 6.6165964726065987
$ws.Range("G4").Value = 5.6686326881838296
$ws.Range("H4").Value = 6.7294661864194607
$ws.Range("I4").Value = 7.2207098269445202
$ws.Range("J4").Value = 7.3191488059459031
$ws.Range("K4").Value = 7.3364889416826751
$ws.Range("L4").Value = 7.7198339498137045
$ws.Range("M4").Value = 7.8258279858854918
$ws.Range("N4").Value = 7.5789073543911334
$ws.Range("O4").Value = 7.4985248229203512
$ws.Range("P4").Value = 7.4141082446031374
$ws.Range("Q4").Value = 7.0384645318913508

# ---------------------------------------------------------------------
# 3. Row 5 ("b) ATMs per 100,000 adults"): same treatment.
# ---------------------------------------------------------------------
$ws.Range("D5").Value = 2.4713943209062914
$ws.Range("E5").Value = 4.8880029305544008
$ws.Range("F5").Value = 8.7206187818873584
$ws.Range("G5").Value = 7.5223561738932325
$ws.Range("H5").Value = 12.187222227373827
$ws.Range("I5").Value = 15.844003577108481
$ws.Range("J5").Value = 20.86084920997822
$ws.Range("K5").Value = 24.989132982250201
$ws.Range("L5").Value = 30.387626630476873
$ws.Range("M5").Value = 31.39992710386154
$ws.Range("N5").Value = 33.570520663807748
$ws.Range("O5").Value = 36.978706525491944
$ws.Range("P5").Value = 39.297069085946042
$ws.Range("Q5").Value = 41.869840292276756

# ---------------------------------------------------------------------
# 4. Row 8 (adult resident population): years 2010-2020 (G8:Q8) were
#    revised upward.
# ---------------------------------------------------------------------
$ws.Range("G8").Value = 3722238
$ws.Range("H8").Value = 3774445
$ws.Range("I8").Value = 3850037
$ws.Range("J8").Value = 3921221
$ws.Range("K8").Value = 3993736
$ws.Range("L8").Value = 4067445
$ws.Range("M8").Value = 4140137
$ws.Range("N8").Value = 4209050
$ws.Range("O8").Value = 4280842
$ws.Range("P8").Value = 4356559
$ws.Range("Q8").Value = 4432785

# ---------------------------------------------------------------------
# 5. Row-height adjustments that shipped with this re-save.
# ---------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 27
$ws.Rows.Item(6).RowHeight = 15.75

# ---------------------------------------------------------------------
# 6. Reset the stored cursor/selection back to A1 (the previous save
#    had left it parked on R13).
# ---------------------------------------------------------------------
$ws.Range("A1").Select()
